$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.190.79'
$ws.Range('D3').Value = '2.806.54'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''360.88'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').Value = '''110.61'
$ws.Range('E6').Value = '  -2.64%  '
$ws.Range('D7').Value = '''0.565'
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '''0.599'
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('D10').Value = '''40.38'
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('D11').Value = '''0.0860'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D13').Value = '''19.57'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').Value = '''7.67'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').Value = '3.238.49'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.799.43'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').Value = '''0.955'
$ws.Range('E17').Value = '  +7.36%  '
$ws.Range('D18').Value = '52.107.24'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').Value = '''7.45'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('D21').Value = '''13.14'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('D22').Value = '0.0₃0990'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '''274.77'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('D24').Value = '''70.53'
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').Value = '''2.77'
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('D26').Value = '''26.84'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = '''10.25'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('E29').Value = '  +5.52%  '
$ws.Range('D30').Value = '''2.16'
$ws.Range('E30').Value = '  -4.15%  '
$ws.Range('D31').Value = '''51.70'
$ws.Range('E31').Value = '  +1.81%  '
$ws.Range('D32').Value = '''0.0464'
$ws.Range('E32').Value = '  +2.03%  '
$ws.Range('D33').Value = '''34.62'
$ws.Range('E33').Value = '  +1.68%  '
$ws.Range('D34').Value = '''5.82'
$ws.Range('E34').Value = '  -0.97%  '
$ws.Range('D35').Value = '''0.0856'
$ws.Range('E35').Value = '  +3.59%  '
$ws.Range('D36').Value = '''5.30'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').Value = '''3.25'
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('D39').Value = '''18.51'
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('E40').Value = '  -2.79%  '
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('D42').Value = '''0.116'
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('D43').Value = '''123.08'
$ws.Range('E43').Value = '  -2.34%  '
$ws.Range('D44').Value = '''2.26'
$ws.Range('E44').Value = '  -2.18%  '
$ws.Range('D45').Value = '''22.27'
$ws.Range('E45').Value = '  -7.14%  '
$ws.Range('D46').Value = '2.085.64'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('D47').Value = '''3.29'
$ws.Range('E47').Value = '  -2.75%  '
$ws.Range('E48').Value = '  -2.21%  '
$ws.Range('D49').Value = '''5.74'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').Value = '''0.944'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('D51').Value = '''8.96'
$ws.Range('E51').Value = '  +0.38%  '
